# "Insert Link on PPT" - turn the GitHub URL text on slide 12 into a
# clickable hyperlink pointing at the same URL.

$p = $ppt.ActivePresentation
$targetUrl = "https://github.com/AniketDohale/Steganography.git"

$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $targetUrl) {
                $tr.ActionSettings.Item(1).Hyperlink.Address = $targetUrl
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not find the Steganography GitHub URL text to hyperlink"
}
